$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.535.13'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.33%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.650.90'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.69%  '

$ws.Range('E4').Value = '  +0.19%  '

$ws.Range('E5').Value = '  +0.23%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '300.45'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.13%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3783'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.40%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3573'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.16%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '50.73'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.07%  '

$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.228'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.19%  '

$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08114'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.86%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.000'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.15%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.16'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.71%  '

$ws.Range('E14').Value = '  -1.46%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.432'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.37%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001206'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.67%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.642.98'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.03%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '97.21'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.02%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06993'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.37%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.796'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.94%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.50'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.29%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.20%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.65'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.09%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.566.95'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.49%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.496'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.90%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.938'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -5.92%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.71%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '152.43'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.05%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.232'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.03%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.26'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.58%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.838.71'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.71%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.015'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.07%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.143'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +5.46%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.041'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -4.54%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02746'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.69%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08743'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.28%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.007'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.13%  '

$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2456'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.46%  '

$ws.Range('E40').Value = '  +4.18%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.06886'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.36%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6943'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.82%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.323'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.09%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.78'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.99%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6462'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.19%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.19%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.274'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.16%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.934'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.57%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07820'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.53%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '128.54'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.90%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.177'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.25%  '
